$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing salt parameters (in-place value fixes) ---
$ws.Cells.Item(75,2).Value = 0.0787        # NaMethaneS: b0 0.787 -> 0.0787
$ws.Cells.Item(86,2).Value = 0.1796        # Me4NEthaneS: b0 0.1792 -> 0.1796
$ws.Cells.Item(88,4).Value = -0.0374       # Bu4NEthaneS: c 0.0036 -> -0.0374

# --- Insert a new row for "pTolueneSA" before the LipTolueneS block ---
$ws.Rows.Item(92).Insert()

$ws.Cells.Item(92,1).Value = "pTolueneSA"
$ws.Cells.Item(92,2).Value = -0.0366
$ws.Cells.Item(92,3).Value = 0.281
$ws.Cells.Item(92,4).Value = 0.0137
$ws.Cells.Item(92,5).Value = 5
$ws.Cells.Item(92,6).Value = 1
$ws.Cells.Item(92,7).Value = -1
$ws.Cells.Item(92,8).Value = 1
$ws.Cells.Item(92,9).Value = 1

# --- Rename compounds that shifted down one row, appending the "S"/"SA" suffix ---
$ws.Cells.Item(94,1).Value = "NapTolueneS"       # was NapToluene
$ws.Cells.Item(95,1).Value = "KpTolueneS"        # was KpToluene
$ws.Cells.Item(96,1).Value = "2,5Me2BenzeneSA"   # was 2,5Me2BenzeneS

# --- Correction to Me3SI parameters (decimal-place fix) ---
$ws.Cells.Item(113,2).Value = -0.0601
$ws.Cells.Item(113,4).Value = 0.0006

# --- Restore the view state (active selection) ---
$ws.Range("C113").Select()
